$d = $word.ActiveDocument

$d.Content.Find.Execute("679÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "394÷3=", 2)
$d.Content.Find.Execute("961÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "878÷9=", 2)
$d.Content.Find.Execute("676÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "181÷9=", 2)
$d.Content.Find.Execute("536÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "419÷8=", 2)
$d.Content.Find.Execute("616÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "543÷7=", 2)
$d.Content.Find.Execute("861÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "226÷4=", 2)
$d.Content.Find.Execute("469÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "897÷5=", 2)
$d.Content.Find.Execute("353÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "994÷6=", 2)
$d.Content.Find.Execute("366÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "447÷6=", 2)
$d.Content.Find.Execute("616÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "547÷5=", 2)
$d.Content.Find.Execute("827÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "195÷3=", 2)
$d.Content.Find.Execute("466÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "769÷5=", 2)
$d.Content.Find.Execute("369÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "701÷8=", 2)
$d.Content.Find.Execute("981÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "367÷7=", 2)
$d.Content.Find.Execute("484÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "411÷4=", 2)
$d.Content.Find.Execute("801÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "820÷7=", 2)
$d.Content.Find.Execute("689÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "581÷3=", 2)
$d.Content.Find.Execute("679÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "441÷2=", 2)
$d.Content.Find.Execute("246÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "668÷2=", 2)
$d.Content.Find.Execute("359÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "343÷9=", 2)
$d.Content.Find.Execute("326÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "845÷8=", 2)
$d.Content.Find.Execute("325÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "649÷9=", 2)
$d.Content.Find.Execute("183÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "978÷3=", 2)
$d.Content.Find.Execute("961÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "915÷6=", 2)
$d.Content.Find.Execute("637÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "104÷3=", 2)
